# Powerpoint writer: consolidate text run nodes.
#
# The canonical edit merges runs of identical formatting inside each
# speaker-notes paragraph (e.g. "Here" / " " / "is" / " " / "a" / " " /
# "note." -> "Here " / "is " / "a " / "note.") so fewer <a:r> nodes are
# written out. The actual rendered text of every notes paragraph is
# unchanged by this refactor - only the run boundaries move - so we
# reconstruct each notes paragraph's full text (paragraph separator is a
# newline) and reassign it on the relevant "Notes Placeholder" shape.

$p = $ppt.ActivePresentation

function Set-NotesBody($slideIndex, $paragraphs) {
    $slide = $p.Slides.Item($slideIndex)
    $notesPage = $slide.NotesPage
    $shape = $notesPage.Shapes.Item(2)
    $shape.TextFrame.TextRange.Text = [string]::Join("`n", $paragraphs)
}

# Slide 1 -> notesSlide1.xml
Set-NotesBody 1 @(
    "Here is a note.",
    "",
    "Here is some other formatting."
)

# Slide 3 -> notesSlide2.xml
Set-NotesBody 3 @(
    "The first note div",
    "",
    "The second note div"
)

# Slide 4 -> notesSlide3.xml
Set-NotesBody 4 @(
    "No link here.",
    "",
    "No note here."
)
